# Apply the "output generated at 456a3b4" data refresh: updates to column F
# ("想去人数" / interested-count) and column G ("最低票价" / lowest price,
# which for one event switched to the text "已停售" = "sales stopped")
# across the 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) sheets.
# 本地生活 (sheet3) has no changes in this update.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 1195
$ws1.Cells.Item(2, 7).Value = "已停售"
$ws1.Cells.Item(3, 7).Value = 68
$ws1.Cells.Item(5, 6).Value = 1289
$ws1.Cells.Item(7, 6).Value = 58
$ws1.Cells.Item(8, 6).Value = 146
$ws1.Cells.Item(9, 6).Value = 360
$ws1.Cells.Item(10, 6).Value = 144
$ws1.Cells.Item(11, 6).Value = 113
$ws1.Cells.Item(12, 6).Value = 895
$ws1.Cells.Item(14, 6).Value = 146
$ws1.Cells.Item(16, 6).Value = 119
$ws1.Cells.Item(17, 6).Value = 350
$ws1.Cells.Item(18, 6).Value = 285
$ws1.Cells.Item(20, 6).Value = 99
$ws1.Cells.Item(21, 6).Value = 680
$ws1.Cells.Item(23, 6).Value = 52
$ws1.Cells.Item(24, 6).Value = 937
$ws1.Cells.Item(25, 6).Value = 388
$ws1.Cells.Item(26, 6).Value = 209
$ws1.Cells.Item(27, 6).Value = 66
$ws1.Cells.Item(30, 6).Value = 26

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(5, 6).Value = 23
$ws2.Cells.Item(7, 6).Value = 268

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 1195
$ws4.Cells.Item(3, 7).Value = "已停售"
$ws4.Cells.Item(4, 7).Value = 68
$ws4.Cells.Item(6, 6).Value = 1289
$ws4.Cells.Item(9, 6).Value = 58
$ws4.Cells.Item(10, 6).Value = 146
$ws4.Cells.Item(11, 6).Value = 360
$ws4.Cells.Item(12, 6).Value = 144
$ws4.Cells.Item(13, 6).Value = 113
$ws4.Cells.Item(14, 6).Value = 895
$ws4.Cells.Item(16, 6).Value = 146
$ws4.Cells.Item(20, 6).Value = 119
$ws4.Cells.Item(21, 6).Value = 23
$ws4.Cells.Item(22, 6).Value = 350
$ws4.Cells.Item(24, 6).Value = 268
$ws4.Cells.Item(25, 6).Value = 285
$ws4.Cells.Item(27, 6).Value = 99
$ws4.Cells.Item(28, 6).Value = 680
$ws4.Cells.Item(30, 6).Value = 52
$ws4.Cells.Item(31, 6).Value = 937
$ws4.Cells.Item(32, 6).Value = 388
$ws4.Cells.Item(35, 6).Value = 209
$ws4.Cells.Item(36, 6).Value = 66
$ws4.Cells.Item(40, 6).Value = 18
$ws4.Cells.Item(41, 6).Value = 26

Write-Output "Edits applied successfully"